# Auto-generated edit script: apply numeric cell updates per the target diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1868.75
$ws.Range("I40").Value = 1754.5454
$ws.Range("K40").Value = 1754.5454
$ws.Range("M40").Value = -1579.5454
$ws.Range("H137").Value = 4878922
$ws.Range("I137").Value = 849.96
$ws.Range("J137").Value = 12500910
$ws.Range("K137").Value = 2549.88
$ws.Range("L137").Value = 37502730
$ws.Range("M137").Value = 0.1199999999998909
$ws.Range("N137").Value = -37507830
$ws.Range("H138").Value = 3578.5356
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3578.5356
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 10735.6068
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -21015.6068

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8894.169
$ws.Range("I32").Value = 8626.207
$ws.Range("J32").Value = 11114.429
$ws.Range("K32").Value = 8626.207
$ws.Range("L32").Value = 11114.429
$ws.Range("M32").Value = -8339.207
$ws.Range("N32").Value = -11688.429
$ws.Range("H45").Value = 3115.5
$ws.Range("I45").Value = 3756
$ws.Range("J45").Value = 2475
$ws.Range("K45").Value = 3756
$ws.Range("L45").Value = 2475
$ws.Range("M45").Value = -3379
$ws.Range("N45").Value = -3229
$ws.Range("H61").Value = 9616694
$ws.Range("I61").Value = 11905964
$ws.Range("J61").Value = 1759.3
$ws.Range("K61").Value = 11905964
$ws.Range("L61").Value = 1759.3
$ws.Range("M61").Value = -11905752
$ws.Range("N61").Value = -2183.3
$ws.Range("H136").Value = 9616694
$ws.Range("I136").Value = 11905964
$ws.Range("J136").Value = 1759.3
$ws.Range("K136").Value = 35717892
$ws.Range("L136").Value = 5277.9
$ws.Range("M136").Value = -35715342
$ws.Range("N136").Value = -10377.9

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4677.037
$ws.Range("I105").Value = 3928
$ws.Range("K105").Value = 3928
$ws.Range("M105").Value = -2181
$ws.Range("H112").Value = 37635.668
$ws.Range("J112").Value = 37635.668
$ws.Range("L112").Value = 37635.668
$ws.Range("N112").Value = -40589.668
$ws.Range("H134").Value = 2646.147
$ws.Range("I134").Value = 1657.5862
$ws.Range("K134").Value = 4972.7586
$ws.Range("M134").Value = -2437.7586

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6064256.5
$ws.Range("I31").Value = 4329.263
$ws.Range("J31").Value = 19609976
$ws.Range("K31").Value = 4329.263
$ws.Range("L31").Value = 19609976
$ws.Range("M31").Value = -4034.263
$ws.Range("N31").Value = -19610566
$ws.Range("H34").Value = 6064256.5
$ws.Range("I34").Value = 4329.263
$ws.Range("J34").Value = 19609976
$ws.Range("K34").Value = 4329.263
$ws.Range("L34").Value = 19609976
$ws.Range("M34").Value = -4127.263
$ws.Range("N34").Value = -19610380
$ws.Range("H107").Value = 563.7083
$ws.Range("I107").Value = 585.5789
$ws.Range("J107").Value = 480.6
$ws.Range("K107").Value = 585.5789
$ws.Range("L107").Value = 480.6
$ws.Range("M107").Value = 1334.4211
$ws.Range("N107").Value = -4320.6
$ws.Range("H140").Value = 44139.89
$ws.Range("J140").Value = 44139.89
$ws.Range("L140").Value = 44139.89
$ws.Range("N140").Value = -54499.89

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 2500
$ws.Range("J58").Value = 2500
$ws.Range("L58").Value = 7500
$ws.Range("N58").Value = -7756
$ws.Range("H113").Value = 2437.7
$ws.Range("I113").Value = 475.8
$ws.Range("J113").Value = 4399.6
$ws.Range("K113").Value = 1427.4
$ws.Range("L113").Value = 13198.8
$ws.Range("M113").Value = 742.5999999999999
$ws.Range("N113").Value = -17538.8
$ws.Range("H131").Value = 857.17
$ws.Range("I131").Value = 725
$ws.Range("J131").Value = 859.8674
$ws.Range("K131").Value = 2175
$ws.Range("L131").Value = 2579.6022
$ws.Range("M131").Value = 2865
$ws.Range("N131").Value = -12659.6022

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1000011
$ws.Range("I113").Value = 1000011
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1000011
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -997841
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 4610.067
$ws.Range("I126").Value = 3381
$ws.Range("K126").Value = 10143
$ws.Range("M126").Value = -7673
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5429.3335
$ws.Range("I7").Value = 5441.25
$ws.Range("J7").Value = 5405.5
$ws.Range("K7").Value = 5441.25
$ws.Range("L7").Value = 5405.5
$ws.Range("M7").Value = -5329.25
$ws.Range("N7").Value = -5629.5
$ws.Range("H40").Value = 7785
$ws.Range("I40").Value = 25850
$ws.Range("J40").Value = 4500.4546
$ws.Range("K40").Value = 25850
$ws.Range("L40").Value = 4500.4546
$ws.Range("M40").Value = -25714
$ws.Range("N40").Value = -4772.4546
$ws.Range("H122").Value = 4970.28
$ws.Range("I122").Value = 5697.2856
$ws.Range("J122").Value = 4045
$ws.Range("K122").Value = 17091.8568
$ws.Range("L122").Value = 12135
$ws.Range("M122").Value = -14641.8568
$ws.Range("N122").Value = -17035
$ws.Range("H126").Value = 5429.3335
$ws.Range("I126").Value = 5441.25
$ws.Range("J126").Value = 5405.5
$ws.Range("K126").Value = 16323.75
$ws.Range("L126").Value = 16216.5
$ws.Range("M126").Value = -13853.75
$ws.Range("N126").Value = -21156.5
$ws.Range("H136").Value = 10644485
$ws.Range("I136").Value = 17859608
$ws.Range("J136").Value = 11673.421
$ws.Range("K136").Value = 53578824
$ws.Range("L136").Value = 35020.263
$ws.Range("M136").Value = -53576274
$ws.Range("N136").Value = -40120.263

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1021.2941
$ws.Range("I107").Value = 1165.8462
$ws.Range("J107").Value = 551.5
$ws.Range("K107").Value = 3497.5386
$ws.Range("L107").Value = 1654.5
$ws.Range("M107").Value = -1577.5386
$ws.Range("N107").Value = -5494.5
$ws.Range("H126").Value = 11157.714
$ws.Range("I126").Value = 6352
$ws.Range("J126").Value = 13080
$ws.Range("K126").Value = 19056
$ws.Range("L126").Value = 39240
$ws.Range("M126").Value = -16586
$ws.Range("N126").Value = -44180
